# Iceland Premier League - base update (16-05-2024)
#
# The underlying dataset was re-sorted; a number of data rows "traded places"
# with each other (their B:AB content moved to a different row), while the
# leftmost sequential index column (A) stayed put. This script reproduces
# that by rotating the B:AB row-content among each affected group of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array lists a group (cycle) of row numbers whose data (columns
# B through AB) rotate among themselves: the content that currently lives in
# row cycle[i+1] ends up in row cycle[i] (wrapping around at the end of the
# group).
$cycles = @(
    @(5,7),
    @(14,15),
    @(21,22),
    @(36,38,37),
    @(72,73,74),
    @(85,86),
    @(91,94),
    @(93,95),
    @(96,97),
    @(98,99),
    @(102,105,103,106),
    @(109,110),
    @(124,125),
    @(133,134),
    @(140,141)
)

foreach ($cycle in $cycles) {
    $n = $cycle.Length

    # Snapshot the current (pre-edit) B:AB values for every row in this group
    # before writing anything, so later writes don't clobber data we still
    # need to read.
    $snapshot = @{}
    foreach ($r in $cycle) {
        $rng = $ws.Range("B$r`:AB$r")
        $snapshot[$r] = $rng.Value2
    }

    # Rotate: new content of cycle[i] = old content of cycle[i+1] (wrap).
    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $cycle[$i]
        $srcRow = $cycle[($i + 1) % $n]
        $ws.Range("B$destRow`:AB$destRow").Value = $snapshot[$srcRow]
    }
}
